$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BG (numeric) and BH (image filename) for rows 1-85, mirroring
# the existing paired numeric/filename column layout already used across
# the sheet (e.g. A/B, C/D, ... AY/AZ, BE/BF). Row 1 additionally gets a
# second new pair in BI/BJ.

$bgValues = @(
    3092.319768143945,
    2892.766574085256,
    1382.381354908828,
    1585.859195749306,
    3071.48895060983,
    705.2288370969691,
    1756.430382804009,
    2316.447288974762,
    774.9667914503081,
    1857.565511411666,
    3432.254558628186,
    1322.304069556817,
    2103.308779132955,
    1784.50670208912,
    3243.2677212897,
    1555.669605120155,
    2628.305760173892,
    2736.082598719962,
    2841.142374109408,
    2453.508030431108,
    2108.441009539911,
    1711.145996860283,
    1949.643762830576,
    2670.571187054704,
    2157.951938171719,
    1914.32194179447,
    1059.654631083203,
    1575.292839029103,
    996.2564907619854,
    1583.444028498974,
    2029.646177997827,
    2214.708368554523,
    714.2857142857144,
    1993.418669242845,
    1809.865958217607,
    1664.65402729139,
    979.9541118222439,
    1676.126071730468,
    665.6804733727812,
    1151.732882502114,
    942.8209153483881,
    1076.862697741819,
    1165.620094191523,
    935.5754135973918,
    1112.788310590509,
    2262.407921748581,
    1936.058447047458,
    1610.61466006521,
    1539.065330274122,
    2127.762347542568,
    3037.978505011473,
    1372.7206859075,
    1259.509721048183,
    2652.457432677214,
    1091.353701243811,
    1691.220867045043,
    1678.843134887091,
    971.8029223523731,
    734.5127400072456,
    2896.691220867046,
    787.34452360826,
    735.7203236324117,
    2321.579519381718,
    1222.678420480619,
    1764.883468180172,
    2736.988286438836,
    2444.451153242363,
    1727.146479893733,
    717.3046733486295,
    899.0460089361191,
    1623.898079942036,
    2033.872720685908,
    1065.390653302741,
    1504.347301050598,
    992.0299480739043,
    930.7450790967276,
    704.9269411906776,
    1319.285110493902,
    922.8957855331483,
    791.5710662963412,
    950.3683130056759,
    1115.203477840841,
    594.1311435816932,
    655.1141166525783,
    1259.811616954474
)

for ($r = 1; $r -le $bgValues.Length; $r++) {
    $ws.Cells.Item($r, 59).Value = $bgValues[$r - 1]
    $ws.Cells.Item($r, 60).Value = "S-CTL4-1_0008.jpg"
}

# Row 1 gets one more pair: BI1 (numeric) / BJ1 (filename).
$ws.Cells.Item(1, 61).Value = 3370.667793744718
$ws.Cells.Item(1, 62).Value = "CK2_2_0006.jpg"
